# Update schedule: shift week of dates forward by one week, update shift
# assignments, fix Madison Johnson's rank, drop Brent Horwitz's old row
# (names shift up by one as Nathan Debergh / Phillip Thompson move up),
# and add two new blank duration/time-formatted cells at the bottom.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: date headers move forward one week ---
$ws.Range("C3").Value = "July9th"
$ws.Range("D3").Value = "July10th"
$ws.Range("E3").Value = "July11th"
$ws.Range("F3").Value = "July12th"
$ws.Range("G3").Value = "July13th"
$ws.Range("H3").Value = "July14th"
$ws.Range("I3").Value = "July15th"

# --- Row 5 (Blake Butz) ---
$ws.Range("E5").Value = "3:30-8"
$ws.Range("F5").Value = "Meet"
$ws.Range("H5").Value = "3:30-8"

# --- Row 6 (Kate North) ---
$ws.Range("C6").Value = "3:30-8"
$ws.Range("E6").Value = "3:30-8"
$ws.Range("F6").ClearContents()
$ws.Range("G6").Value = "X"
$ws.Range("H6").Value = "3:30-8"

# --- Row 7 (Emerson Metzger) ---
$ws.Range("C7").Value = "X"
$ws.Range("D7").Value = "X"
$ws.Range("E7").ClearContents()
$ws.Range("F7").ClearContents()
$ws.Range("G7").ClearContents()
$ws.Range("H7").ClearContents()

# --- Row 8 (Avery Larsen) ---
$ws.Range("C8").Value = "10:30-3:30"
$ws.Range("D8").Value = "10:30-3:30"
$ws.Range("E8").Value = "3:30-8"
$ws.Range("F8").Value = "Meet"
$ws.Range("G8").Value = "X"
$ws.Range("H8").Value = "X"

# --- Row 9 (Austin Page) ---
$ws.Range("C9").Value = "10:30-3:30"
$ws.Range("D9").Value = "3:30-8"
$ws.Range("E9").Value = "OFF"
$ws.Range("G9").Value = "OFF"
$ws.Range("H9").Value = "OFF"

# --- Row 10 (Riley White) ---
$ws.Range("D10").Value = "X"
$ws.Range("E10").Value = "X"
$ws.Range("F10").Value = "Meet"

# --- Row 11 (Robert Wade) ---
$ws.Range("D11").Value = "10:30-3:30"
$ws.Range("E11").ClearContents()
$ws.Range("F11").ClearContents()
$ws.Range("H11").Value = "10:30-3:30"

# --- Row 12 (Tatum Plunk) ---
$ws.Range("C12").Value = "X"
$ws.Range("D12").ClearContents()
$ws.Range("E12").ClearContents()
$ws.Range("F12").ClearContents()

# --- Row 13 (Michael Vangruber) ---
$ws.Range("D13").ClearContents()
$ws.Range("E13").ClearContents()
$ws.Range("F13").Value = "10:30-3:30"

# --- Row 14 (Jackson Blakely) ---
$ws.Range("E14").ClearContents()

# --- Row 15 (Addison Clark) ---
$ws.Range("C15").Value = "X"
$ws.Range("E15").ClearContents()
$ws.Range("F15").Value = "10:30-3:30"
$ws.Range("G15").Value = "10:30-3:30"

# --- Row 16: Nathan Debergh moves up into this row ---
$ws.Range("A16").Value = "Nathan Debergh"
$ws.Range("C16").Value = "X"
$ws.Range("D16").Value = "X"
$ws.Range("E16").Value = "X"
$ws.Range("F16").Value = "X"
$ws.Range("G16").Value = "X"
$ws.Range("H16").Value = "X"

# --- Row 17: Phillip Thompson moves up into this row ---
$ws.Range("A17").Value = "Phillip Thompson"
$ws.Range("C17").ClearContents()
$ws.Range("D17").ClearContents()
$ws.Range("E17").ClearContents()
$ws.Range("F17").ClearContents()
$ws.Range("G17").ClearContents()
$ws.Range("H17").ClearContents()

# --- Row 18: Madison Johnson moves up, rank corrected to 3.0 ---
$ws.Range("A18").Value = "Madison Johnson"
$ws.Range("B18").Value = 3.0
$ws.Range("C18").Value = "X"

# --- Row 19: Asher Bobbett moves up ---
$ws.Range("A19").Value = "Asher Bobbett"
$ws.Range("C19").Value = "10:30-3:30"
$ws.Range("D19").ClearContents()
$ws.Range("E19").ClearContents()
$ws.Range("F19").ClearContents()
$ws.Range("H19").Value = "X"

# --- Row 20: Blake Ucherek moves up ---
$ws.Range("A20").Value = "Blake Ucherek"
$ws.Range("D20").ClearContents()
$ws.Range("E20").ClearContents()
$ws.Range("F20").Value = "10:30-3:30"
$ws.Range("G20").ClearContents()
$ws.Range("H20").ClearContents()

# --- Row 21: Ethan Van Horn moves up ---
$ws.Range("A21").Value = "Ethan Van Horn "
$ws.Range("E21").ClearContents()

# --- Row 22: Kai King moves up ---
$ws.Range("A22").Value = "Kai King"
$ws.Range("D22").ClearContents()
$ws.Range("E22").ClearContents()
$ws.Range("F22").Value = "Meet"

# --- Row 23: Madeline Ellison moves up ---
$ws.Range("A23").Value = "Madeline Ellison"
$ws.Range("E23").ClearContents()

# --- Row 24: Tyler Carpenter moves up ---
$ws.Range("A24").Value = "Tyler Carpenter"
$ws.Range("C24").ClearContents()
$ws.Range("D24").ClearContents()
$ws.Range("E24").ClearContents()
$ws.Range("F24").ClearContents()

# --- Row 25 (Jayden Garcia) stays but schedule updates ---
$ws.Range("C25").Value = "10:30-3:30"
$ws.Range("E25").Value = "10:30-3:30"
$ws.Range("F25").Value = "Meet"
$ws.Range("G25").Value = "10:30-3:30"

# --- Row 26 (Naya Okonkwo) stays ---
$ws.Range("E26").ClearContents()

# --- Row 27 (Bella Hamilton) stays ---
$ws.Range("C27").Value = "3:30-8"
$ws.Range("D27").ClearContents()
$ws.Range("E27").ClearContents()

# --- Row 28: Brent Horwitz moves up ---
$ws.Range("A28").Value = "Brent Horwitz"
$ws.Range("E28").ClearContents()
$ws.Range("G28").Value = "X"
$ws.Range("H28").Value = "X"

# --- Row 29 removed entirely (old Brent Horwitz row) ---
$ws.Rows.Item(29).Delete()

# --- New rows 33 & 34: blank cells with duration / time number formats ---
$ws.Range("A33").NumberFormat = "[h]:mm:ss"
$ws.Range("C33").NumberFormat = "[h]:mm:ss"
$ws.Range("A34").NumberFormat = "[h]:mm:ss"
$ws.Range("B34").NumberFormat = "h:mm"
